# Auto-generated edit script applying numeric corrections to the
# Marilith_Profits leve-profit tables, per sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 170.14285
$ws.Range("I6").Value = 114
$ws.Range("J6").Value = 245
$ws.Range("K6").Value = 342
$ws.Range("L6").Value = 735
$ws.Range("M6").Value = -230
$ws.Range("N6").Value = -959
$ws.Range("H28").Value = 3599.7144
$ws.Range("I28").Value = 899.6667
$ws.Range("K28").Value = 899.6667
$ws.Range("M28").Value = -414.6667
$ws.Range("H40").Value = 2385.5715
$ws.Range("H48").Value = 4999
$ws.Range("J48").Value = 4999
$ws.Range("L48").Value = 14997
$ws.Range("N48").Value = -15581
$ws.Range("H55").Value = 375
$ws.Range("I55").Value = 450
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 450
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = -236
$ws.Range("N55").Value = -728
$ws.Range("H56").Value = 4999
$ws.Range("J56").Value = 4999
$ws.Range("L56").Value = 14997
$ws.Range("N56").Value = -16065
$ws.Range("H62").Value = 3924
$ws.Range("I62").Value = 2565.6667
$ws.Range("K62").Value = 2565.6667
$ws.Range("M62").Value = -1941.6667
$ws.Range("H65").Value = 3924
$ws.Range("I65").Value = 2565.6667
$ws.Range("K65").Value = 12828.3335
$ws.Range("M65").Value = -9708.333500000001
$ws.Range("H70").Value = 2813.9
$ws.Range("I70").Value = 1336.875
$ws.Range("J70").Value = 3798.5833
$ws.Range("K70").Value = 4010.625
$ws.Range("L70").Value = 11395.7499
$ws.Range("M70").Value = -3740.625
$ws.Range("N70").Value = -11935.7499
$ws.Range("H73").Value = 2813.9
$ws.Range("I73").Value = 1336.875
$ws.Range("J73").Value = 3798.5833
$ws.Range("K73").Value = 4010.625
$ws.Range("L73").Value = 11395.7499
$ws.Range("M73").Value = -3074.625
$ws.Range("N73").Value = -13267.7499
$ws.Range("H86").Value = 13988.556
$ws.Range("I86").Value = 8899.5
$ws.Range("J86").Value = 24166.666
$ws.Range("K86").Value = 8899.5
$ws.Range("L86").Value = 24166.666
$ws.Range("M86").Value = -7776.5
$ws.Range("N86").Value = -26412.666
$ws.Range("H89").Value = 13988.556
$ws.Range("I89").Value = 8899.5
$ws.Range("J89").Value = 24166.666
$ws.Range("K89").Value = 44497.5
$ws.Range("L89").Value = 120833.33
$ws.Range("M89").Value = -38881.5
$ws.Range("N89").Value = -132065.33
$ws.Range("H132").Value = 2499.3076
$ws.Range("I132").Value = 2591
$ws.Range("J132").Value = 1399
$ws.Range("K132").Value = 7773
$ws.Range("L132").Value = 4197
$ws.Range("M132").Value = -5243
$ws.Range("N132").Value = -9257
$ws.Range("H137").Value = 2983.1428
$ws.Range("I137").Value = 1964.6666
$ws.Range("K137").Value = 5893.9998
$ws.Range("M137").Value = -3343.9998
$ws.Range("H138").Value = 2749.25
$ws.Range("I138").Value = 1831.8334
$ws.Range("J138").Value = 3666.6667
$ws.Range("K138").Value = 5495.5002
$ws.Range("L138").Value = 11000.0001
$ws.Range("M138").Value = -355.5002000000004
$ws.Range("N138").Value = -21280.0001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H88").Value = 2826.2144
$ws.Range("I88").Value = 1128
$ws.Range("J88").Value = 4099.875
$ws.Range("K88").Value = 1128
$ws.Range("L88").Value = 4099.875
$ws.Range("M88").Value = -722
$ws.Range("N88").Value = -4911.875
$ws.Range("H91").Value = 2826.2144
$ws.Range("I91").Value = 1128
$ws.Range("J91").Value = 4099.875
$ws.Range("K91").Value = 1128
$ws.Range("L91").Value = 4099.875
$ws.Range("M91").Value = 276
$ws.Range("N91").Value = -6907.875
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 44699.75
$ws.Range("J130").Value = 44699.75
$ws.Range("L130").Value = 44699.75
$ws.Range("N130").Value = -54739.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 30249.5
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 30249.5
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 30249.5
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -30589.5
$ws.Range("H20").Value = 51000
$ws.Range("J20").Value = 51000
$ws.Range("L20").Value = 51000
$ws.Range("N20").Value = -51472
$ws.Range("H30").Value = 51000
$ws.Range("J30").Value = 51000
$ws.Range("L30").Value = 51000
$ws.Range("N30").Value = -51182
$ws.Range("H31").Value = 3504.6086
$ws.Range("I31").Value = 2650.1765
$ws.Range("K31").Value = 2650.1765
$ws.Range("M31").Value = -2355.1765
$ws.Range("H34").Value = 3504.6086
$ws.Range("I34").Value = 2650.1765
$ws.Range("K34").Value = 2650.1765
$ws.Range("M34").Value = -2448.1765
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H128").Value = 51000
$ws.Range("J128").Value = 51000
$ws.Range("L128").Value = 51000
$ws.Range("N128").Value = -60960
$ws.Range("H134").Value = 3109.2
$ws.Range("I134").Value = 2978.5
$ws.Range("J134").Value = 3632
$ws.Range("K134").Value = 8935.5
$ws.Range("L134").Value = 10896
$ws.Range("M134").Value = -6400.5
$ws.Range("N134").Value = -15966

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 225.1579
$ws.Range("I2").Value = 157.6923
$ws.Range("J2").Value = 371.33334
$ws.Range("K2").Value = 946.1537999999999
$ws.Range("L2").Value = 2228.00004
$ws.Range("M2").Value = -833.1537999999999
$ws.Range("N2").Value = -2454.00004
$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 10000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -35242
$ws.Range("H120").Value = 6138.8335
$ws.Range("I120").Value = 900
$ws.Range("J120").Value = 8758.25
$ws.Range("K120").Value = 2700
$ws.Range("L120").Value = 26274.75
$ws.Range("M120").Value = 2138
$ws.Range("N120").Value = -35950.75

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 50000
$ws.Range("J32").Value = 50000
$ws.Range("L32").Value = 50000
$ws.Range("N32").Value = -50592
$ws.Range("H70").Value = 83336890
$ws.Range("I70").Value = 111114110
$ws.Range("K70").Value = 111114110
$ws.Range("M70").Value = -111113840
$ws.Range("H73").Value = 83336890
$ws.Range("I73").Value = 111114110
$ws.Range("K73").Value = 111114110
$ws.Range("M73").Value = -111113174

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 6690.5454
$ws.Range("J132").Value = 8400
$ws.Range("L132").Value = 25200
$ws.Range("N132").Value = -30260
$ws.Range("H136").Value = 3251.8333
$ws.Range("I136").Value = 3252.75
$ws.Range("K136").Value = 9758.25
$ws.Range("M136").Value = -7208.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2100
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 6000
$ws.Range("N100").Value = -7082
$ws.Range("H122").Value = 2795.6
$ws.Range("I122").Value = 2499.5
$ws.Range("J122").Value = 2993
$ws.Range("K122").Value = 7498.5
$ws.Range("L122").Value = 8979
$ws.Range("M122").Value = -5048.5
$ws.Range("N122").Value = -13879
$ws.Range("H126").Value = 2049.25
$ws.Range("I126").Value = 2049.25
$ws.Range("K126").Value = 6147.75
$ws.Range("M126").Value = -3677.75
$ws.Range("H132").Value = 1326.4814
$ws.Range("I132").Value = 1127.2174
$ws.Range("K132").Value = 3381.6522
$ws.Range("M132").Value = -851.6522

